$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" date placeholders (slide master +
#    every slide layout) from 5/24/2024 -> 5/28/2024, mirroring what
#    PowerPoint itself writes back into every layout/master on save.
# ---------------------------------------------------------------------------
function Set-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "5/28/2024"
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholders($master.Shapes)
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholders($layout.Shapes)
}

# ---------------------------------------------------------------------------
# 2. Append two new "Title and Content" slides (layout 16) at the end of the
#    deck, matching the new slide10.xml / slide11.xml content.
# ---------------------------------------------------------------------------

# --- Slide 10: "A Description of the Code" ---
$s10 = $p.Slides.Add($p.Slides.Count + 1, 16)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "A Description of the Code"

$body10 = $s10.Shapes.Item(2).TextFrame.TextRange
$body10.Text = "This program looks at the world’s worst ten carbon producers.  You can specify which of these countries you’d like a look at:`r" + `
    "country_of_interest = input (‘Which country are you interested in?’)`r" + `
    "Earlier, the dataframe is cleaned:`r" + `
    "df.dropna (axis = 0, how = ‘any’, subset = ‘Country’, inplace = True)`r" + `
    "df [‘Country’] = df [‘Country’].str.replace (r’\W’, ‘’, regex = True)`r" + `
    "and a bar chart of the top ten coal burners in 2021 produced using Seaborn and MatPlotLib:`r" + `
    "df1 = pd.DataFrame ({‘Country’: top_ten, ‘Emissions from Coal Burning (MtCO2)’: top_ten_Coal})`r" + `
    "sns.barplot (x = ‘Country’, y = ‘Emissions from Coal Burning (MtCO2’, data = df1).set (title = ‘Top Ten Coal Burners in 2021’)`r" + `
    "plt.xticks (rotation = 45)`r" + `
    "plt.show ()"

# --- Slide 11: (blank title) "Once you've input which country..." ---
$s11 = $p.Slides.Add($p.Slides.Count + 1, 16)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = ""

$body11 = $s11.Shapes.Item(2).TextFrame.TextRange
$body11.Text = "Once you’ve input which country you’re interested in, a line plot of emissions from burning coal in 2021 for that country is produced:`r" + `
    "df1 = pd.DataFrame ({‘Year’: Country_Year, ‘Emission from Burning Coal (MtCO2)’: Country_Coal})`r" + `
    "sns.lineplot (x = ‘Year’, y = ‘Emission from Burning Coal (MtCO2)’, data = df1).set (title = country_of_interest)`r" + `
    "plt.show ()."

Write-Output ("Slide count now: " + $p.Slides.Count)
